$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '28.535.71'
$ws.Range('E2').Value = '  +5.08%  '

# Row 3
$ws.Range('D3').Value = '1.603.64'
$ws.Range('E3').Value = '  +2.71%  '

# Row 4
$ws.Range('E4').Value = '  -0.36%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.44'
$ws.Range('E5').Value = '  +2.36%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.500'
$ws.Range('E6').Value = '  +1.94%  '

# Row 7
$ws.Range('E7').Value = '  -0.40%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '24.12'
$ws.Range('E8').Value = '  +9.60%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.253'
$ws.Range('E9').Value = '  +1.88%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0603'
$ws.Range('E10').Value = '  +1.11%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0890'
$ws.Range('E11').Value = '  +2.17%  '

# Row 12
$ws.Range('D12').Value = '1.831.90'
$ws.Range('E12').Value = '  +2.70%  '

# Row 13
$ws.Range('D13').Value = '1.608.34'
$ws.Range('E13').Value = '  +2.61%  '

# Row 14
$ws.Range('E14').Value = '  +1.12%  '

# Row 15
$ws.Range('E15').Value = '  +3.56%  '

# Row 16
$ws.Range('D16').Value = '28.533.76'
$ws.Range('E16').Value = '  +5.14%  '

# Row 17
$ws.Range('E17').Value = '  +2.78%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '233.24'
$ws.Range('E18').Value = '  +7.74%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.58'
$ws.Range('E19').Value = '  +1.92%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0712'
$ws.Range('E20').Value = '  +1.56%  '

# Row 21
$ws.Range('E21').Value = '  -0.43%  '

# Row 22
$ws.Range('E22').Value = '  +0.57%  '

# Row 23
$ws.Range('E23').Value = '  +2.55%  '

# Row 24
$ws.Range('E24').Value = '  +1.52%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '152.74'
$ws.Range('E25').Value = '  -0.44%  '

# Row 26
$ws.Range('E26').Value = '  +2.19%  '

# Row 27
$ws.Range('E27').Value = '  +0.21%  '

# Row 28
$ws.Range('E28').Value = '  +1.23%  '

# Row 29
$ws.Range('E29').Value = '  -0.50%  '

# Row 30
$ws.Range('E30').Value = '  +0.90%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0477'
$ws.Range('E31').Value = '  +1.37%  '

# Row 32
$ws.Range('E32').Value = '  +0.79%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.18'
$ws.Range('E33').Value = '  +0.95%  '

# Row 34
$ws.Range('D34').Value = '1.426.30'
$ws.Range('E34').Value = '  -0.54%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.60'
$ws.Range('E35').Value = '  -0.29%  '

# Row 36
$ws.Range('E36').Value = '  -4.26%  '

# Row 37
$ws.Range('E37').Value = '  -0.28%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0168'
$ws.Range('E38').Value = '  +1.17%  '

# Row 39
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.546'
$ws.Range('E39').Value = '  +2.79%  '

# Row 40
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.53'
$ws.Range('E40').Value = '  +8.21%  '

# Row 41
$ws.Range('E41').Value = '  +2.39%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.75'
$ws.Range('E42').Value = '  -2.83%  '

# Row 43
$ws.Range('E43').Value = '  -0.52%  '

# Row 44
$ws.Range('E44').Value = '  +6.95%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.978'
$ws.Range('E45').Value = '  -2.10%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '65.19'
$ws.Range('E46').Value = '  +1.45%  '

# Row 47
$ws.Range('D47').Value = '1.743.15'
$ws.Range('E47').Value = '  +2.80%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '87.75'
$ws.Range('E48').Value = '  +2.68%  '

# Row 49
$ws.Range('E49').Value = '  +0.24%  '

# Row 50
$ws.Range('D50').Value = '0.0₆0107'
$ws.Range('E50').Value = '  +6.53%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0527'
$ws.Range('E51').Value = '  +0.64%  '
